# Derek's Log.xlsx - add Tuesday (Sep 6) and Wednesday (Sep 7) log entries
# and extend the data validation ranges to cover the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# ---------------------------------------------------------------------
# 1. Give row 115 the slightly different row-height / descent that the
#    new sheet extent ends up with (cosmetic - matches author's resave).
# ---------------------------------------------------------------------
$ws.Rows.Item(115).RowHeight = 14.45

# ---------------------------------------------------------------------
# 2. Stamp out the new rows 116-133 by copying the formatting of an
#    existing row with an identical style signature, then overwrite the
#    cell values/content for the new entries.
# ---------------------------------------------------------------------

# Row 116: Tuesday header/separator row (style like row 112)
$ws.Range("A112:F112").Copy($ws.Range("A116:F116"))
$ws.Range("C116").Value = "TUESDAY"

# Row 117
$ws.Range("A52:F52").Copy($ws.Range("A117:F117"))
$ws.Range("A117").Value = "Other"
$ws.Range("B117").Value = 42619
$ws.Range("C117").Value = "1630"
$ws.Range("D117").Value = "MC"
$ws.Range("E117").Value = "215"
$ws.Range("F117").Value = "Drawer does not lock."

# Row 118
$ws.Range("A3:F3").Copy($ws.Range("A118:F118"))
$ws.Range("A118").Value = "Pickup Mic"
$ws.Range("B118").Value = 42619
$ws.Range("C118").Value = "1600"
$ws.Range("D118").Value = "MC"
$ws.Range("E118").Value = "140-SCR"
$ws.Range("F118").Value = "Door code 7083*  Neck mic and small PA, to Founders 156A."
$ws.Rows.Item(118).RowHeight = 15

# Row 119
$ws.Range("A3:F3").Copy($ws.Range("A119:F119"))
$ws.Range("A119").Value = "Pickup PC"
$ws.Range("B119").Value = 42619
$ws.Range("C119").Value = "1600"
$ws.Range("D119").Value = "MC"
$ws.Range("E119").Value = "140-SCR"
$ws.Range("F119").Value = "Door code 7083*  Leave portable screen. Return to Founders 156A."
$ws.Rows.Item(119).RowHeight = 15

# Row 120
$ws.Range("A3:F3").Copy($ws.Range("A120:F120"))
$ws.Range("A120").Value = "AV Shutdown"
$ws.Range("B120").Value = 42619
$ws.Range("C120").Value = "2100"
$ws.Range("D120").Value = "VC"
$ws.Range("E120").Value = "001-JCR"
$ws.Range("F120").Value = " Leave all in an lock room. Do not pick up equipment. Needed again tomorrow morning."
$ws.Rows.Item(120).RowHeight = 30

# Row 121: Wednesday header/separator row (style like row 110)
$ws.Range("A110:F110").Copy($ws.Range("A121:F121"))
$ws.Range("C121").Value = "WEDNESDAY"

# Row 122
$ws.Range("A3:F3").Copy($ws.Range("A122:F122"))
$ws.Range("A122").Value = "AV Shutdown"
$ws.Range("B122").Value = 42620
$ws.Range("C122").Value = "1600"
$ws.Range("D122").Value = "R"
$ws.Range("E122").Value = "N102"
$ws.Range("F122").Value = "Nat Taylor Cinema. Lock cinema all doors after shutdown."

# Row 123
$ws.Range("A3:F3").Copy($ws.Range("A123:F123"))
$ws.Range("A123").Value = "Pickup Mic"
$ws.Range("B123").Value = 42620
$ws.Range("C123").Value = "1700"
$ws.Range("D123").Value = "SLH"
$ws.Range("E123").Value = "A"
$ws.Range("F123").Value = "Lec mic, stand and cable. To  Stedman 114L MCR."

# Row 124
$ws.Range("A3:F3").Copy($ws.Range("A124:F124"))
$ws.Range("A124").Value = "Pickup Mic"
$ws.Range("B124").Value = 42620
$ws.Range("C124").Value = "1700"
$ws.Range("D124").Value = "SLH"
$ws.Range("E124").Value = "B"
$ws.Range("F124").Value = "Lec mic, stand and cable. To  Stedman 114L MCR."

# Row 125
$ws.Range("A3:F3").Copy($ws.Range("A125:F125"))
$ws.Range("A125").Value = "Pickup Mic"
$ws.Range("B125").Value = 42620
$ws.Range("C125").Value = "1700"
$ws.Range("D125").Value = "SLH"
$ws.Range("E125").Value = "D"
$ws.Range("F125").Value = "Lec mic, stand and cable. To  Stedman 114L MCR."

# Row 126
$ws.Range("A3:F3").Copy($ws.Range("A126:F126"))
$ws.Range("A126").Value = "Pickup PC"
$ws.Range("B126").Value = 42620
$ws.Range("C126").Value = "1700"
$ws.Range("D126").Value = "VC"
$ws.Range("E126").Value = "001-JCR"
$ws.Range("F126").Value = "May include portable screen. All to Vanier 040."
$ws.Rows.Item(126).RowHeight = 15

# Row 127
$ws.Range("A3:F3").Copy($ws.Range("A127:F127"))
$ws.Range("A127").Value = "Pickup Small PA"
$ws.Range("B127").Value = 42620
$ws.Range("C127").Value = "1700"
$ws.Range("D127").Value = "VC"
$ws.Range("E127").Value = "001-JCR"
$ws.Range("F127").Value = "Lectern mic with small PA etc. Return to Vanier 040."

# Row 128
$ws.Range("A3:F3").Copy($ws.Range("A128:F128"))
$ws.Range("A128").Value = "Pickup PC"
$ws.Range("B128").Value = 42620
$ws.Range("C128").Value = "1700"
$ws.Range("D128").Value = "VC"
$ws.Range("E128").Value = "113"
$ws.Range("F128").Value = "Flat screen DLP and wireless keyboard. To Vanier 132 storeroom."

# Row 129
$ws.Range("A3:F3").Copy($ws.Range("A129:F129"))
$ws.Range("A129").Value = "Pickup Mic"
$ws.Range("B129").Value = 42620
$ws.Range("C129").Value = "1700"
$ws.Range("D129").Value = "VC"
$ws.Range("E129").Value = "135"
$ws.Range("F129").Value = "Lec mic, stand and cable. To Vanier 040."

# Row 130
$ws.Range("A3:F3").Copy($ws.Range("A130:F130"))
$ws.Range("A130").Value = "Pickup Mic"
$ws.Range("B130").Value = 42620
$ws.Range("C130").Value = "1700"
$ws.Range("D130").Value = "VC"
$ws.Range("E130").Value = "258"
$ws.Range("F130").Value = "Lec mic, stand and cable, small PA. To  Vanier 040."

# Row 131
$ws.Range("A111:F111").Copy($ws.Range("A131:F131"))
$ws.Range("A131").Value = "Other"
$ws.Range("B131").Value = 42620
$ws.Range("C131").Value = "1730"
$ws.Range("D131").Value = "MC"
$ws.Range("E131").Value = "157A"
$ws.Range("F131").Value = "Door code 11012* "

# Row 132 (rich text comment cell)
$ws.Range("A99:F99").Copy($ws.Range("A132:F132"))
$ws.Range("A132").Value = "Pickup Mic"
$ws.Range("B132").Value = 42620
$ws.Range("C132").Value = "2045"
$ws.Range("D132").Value = "FC"
$ws.Range("E132").Value = "152"
$f132 = $ws.Range("F132")
$f132.Value = "Leave mic cables in place. Remove lec mic, stand, clip. Turn off PA as usual and lock room. To Founders 156A."
$f132.Characters(1, 26).Font.Underline = 2
$f132.Characters(27, 85).Font.Underline = $false
$ws.Rows.Item(132).RowHeight = 30

# Row 133
$ws.Range("A3:F3").Copy($ws.Range("A133:F133"))
$ws.Range("A133").Value = "AV Shutdown"
$ws.Range("B133").Value = 42620
$ws.Range("C133").Value = "2045"
$ws.Range("D133").Value = "FC"
$ws.Range("E133").Value = "152"
$ws.Range("F133").Value = "Return wireless keyboard & projector remote to Founders 156A."

# ---------------------------------------------------------------------
# 3. Extend the three data-validation ranges so the new rows keep the
#    same drop-down lists as the rows above them.
# ---------------------------------------------------------------------
$ws.Range("A111:A1048576").Validation.Delete()
$ws.Range("A111:A120").Validation.Add(3, 1, 1, "=Task_type")
$ws.Range("B121").Validation.Add(3, 1, 1, "=Task_type")
$ws.Range("A122:A1048576").Validation.Add(3, 1, 1, "=Task_type")

$ws.Range("D111:D1048576").Validation.Delete()
$ws.Range("D111:D120").Validation.Add(3, 1, 1, "=Building")
$ws.Range("E121").Validation.Add(3, 1, 1, "=Building")
$ws.Range("D122:D1048576").Validation.Add(3, 1, 1, "=Building")

$ws.Range("A121").Validation.Add(3, 1, 1, "=Staff_Name")

# ---------------------------------------------------------------------
# 4. Update the frozen-pane top row and final selection to match where
#    the author ended up after entering the new data.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("A110").Select()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("C150").Select()
